$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Assignee" (E) column for the existing "Data Collect" block
# (rows 11-13 already hold D = "Data Collect"/blank task-type cells; the new
# data adds the task description + assignee names that were missing).
$ws.Range("D11").Value = "scrap data from internet"
$ws.Range("E11").Value = "Zikri"

$ws.Range("D12").Value = "scrap data from internet"
$ws.Range("E12").Value = "Afif"

$ws.Range("D13").Value = "scrap data from internet"
$ws.Range("E13").Value = "Naufal"

# Row 13 shrinks back down to the "short" row height used by the rest of
# the "scrap data from internet" block (was the tall 44.25 row before).
$ws.Rows.Item(13).RowHeight = 29.25

# --- Append two more repeats of the same 5-row block (rows 14-20, 21-25),
# matching the same Task/Assignee/height cycle already established by
# rows 11-15 (Zikri/Afif/Naufal/Zahran/Satrio).
$names = @("Zikri", "Afif", "Naufal", "Zahran", "Satrio")
$heights = @(30.75, 29.25, 29.25, 44.25, 44.25)

for ($r = 14; $r -le 25; $r++) {
    $i = ($r - 11) % 5

    # Clone formatting (borders/fill/alignment => style index 1) from row 11.
    $src = $ws.Range("A11:H11")
    $dst = $ws.Range("A$r:H$r")
    $src.Copy($dst)

    $ws.Rows.Item($r).RowHeight = $heights[$i]

    $ws.Range("D$r").Value = "scrap data from internet"
    $ws.Range("E$r").Value = $names[$i]
}

# --- Restore the selection/scroll state recorded in the saved file.
$ws.Range("A21:XFD25").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
